$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between row 109 (Twente-PSV) and
#     row 110 (Heerenveen-Sittard). Columns A-E (index/pais/torneio/
#     temporada/data_partida) stay put on each row. ---
$row109 = $ws.Range("F109:V109").Value2
$row110 = $ws.Range("F110:V110").Value2
$ws.Range("F109:V109").Value2 = $row110
$ws.Range("F110:V110").Value2 = $row109

# --- Append the new match row 117 (Heerenveen 3-0 Almere City) ---
$ws.Range("A117").Value2 = 116
$ws.Range("B117").Value2 = "netherlands"
$ws.Range("C117").Value2 = "eredivisie"
$ws.Range("D117").Value2 = "2023-2024"
$ws.Range("E117").Value2 = 45261.83333333334
$ws.Range("F117").Value2 = "Heerenveen"
$ws.Range("G117").Value2 = 3
$ws.Range("H117").Value2 = "Almere City"
$ws.Range("I117").Value2 = 0
$ws.Range("J117").Value2 = 1.6
$ws.Range("K117").Value2 = "26/11/2023 12:42"
$ws.Range("L117").Value2 = 1.59
$ws.Range("M117").Value2 = "01/12/2023 19:55"
$ws.Range("N117").Value2 = 4.39
$ws.Range("O117").Value2 = "26/11/2023 12:42"
$ws.Range("P117").Value2 = 4.24
$ws.Range("Q117").Value2 = "01/12/2023 19:57"
$ws.Range("R117").Value2 = 5.51
$ws.Range("S117").Value2 = "26/11/2023 12:42"
$ws.Range("T117").Value2 = 6.01
$ws.Range("U117").Value2 = "01/12/2023 19:57"
$ws.Range("V117").Value2 = "https://www.betexplorer.com/football/netherlands/eredivisie/heerenveen-almere-city/tpaayL36/"

# Copy the formatting used by the rest of the table (bold/bordered index
# cell in column A, date/time number format in column E) onto the new
# row, matching the pattern of every other data row.
$ws.Range("A116").Copy()
$ws.Range("A117").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E116").Copy()
$ws.Range("E117").PasteSpecial(-4122)  # xlPasteFormats
